$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.404689
$ws.Range("H2").Value = 16.214067
$ws.Range("I2").Value = 0.08747859311663772
$ws.Range("J2").Value = 0.09021076583983562
$ws.Range("M2").Value = 0.4260053333333333
$ws.Range("N2").Value = 1.278016
$ws.Range("O2").Value = 0.02405532912416773
$ws.Range("P2").Value = 0.02531756756689831
$ws.Range("Q2").Value = 2.302426339008
$ws.Range("R2").Value = 20.721837051072
$ws.Range("S2").Value = 0.002104326348739873
$ws.Range("T2").Value = 0.00228391715941168
$ws.Range("G3").Value = 5.404689
$ws.Range("H3").Value = 16.214067
$ws.Range("I3").Value = 0.08747859311663772
$ws.Range("J3").Value = 0.09021076583983562
$ws.Range("O3").Value = 0.05879323641880037
$ws.Range("P3").Value = 0.06187825274916518
$ws.Range("Q3").Value = 5.627322552414
$ws.Range("R3").Value = 50.645902971726
$ws.Range("S3").Value = 0.005143149606690524
$ws.Range("T3").Value = 0.005582084569333105
$ws.Range("G4").Value = 5.404689
$ws.Range("H4").Value = 16.214067
$ws.Range("I4").Value = 0.08747859311663772
$ws.Range("J4").Value = 0.09021076583983562
$ws.Range("M4").Value = 5.850740666666667
$ws.Range("N4").Value = 17.552222
$ws.Range("O4").Value = 0.3303749538898241
$ws.Range("P4").Value = 0.3477104875323931
$ws.Range("Q4").Value = 31.621433722986
$ws.Range("R4").Value = 284.592903506874
$ws.Range("S4").Value = 0.02890073616725587
$ws.Range("T4").Value = 0.0313672293708398
$ws.Range("G5").Value = 5.404689
$ws.Range("H5").Value = 16.214067
$ws.Range("I5").Value = 0.08747859311663772
$ws.Range("J5").Value = 0.09021076583983562
$ws.Range("M5").Value = 2.648771
$ws.Range("N5").Value = 5.297542
$ws.Range("O5").Value = 0.1495686865725097
$ws.Range("P5").Value = 0.1049445996947469
$ws.Range("Q5").Value = 14.315783487219
$ws.Range("R5").Value = 85.89470092331401
$ws.Range("S5").Value = 0.0130840582756665
$ws.Range("T5").Value = 0.009467132709218094
$ws.Range("G6").Value = 5.404689
$ws.Range("H6").Value = 16.214067
$ws.Range("I6").Value = 0.08747859311663772
$ws.Range("J6").Value = 0.09021076583983562
$ws.Range("M6").Value = 7.742685666666667
$ws.Range("N6").Value = 23.228057
$ws.Range("O6").Value = 0.4372077939946981
$ws.Range("P6").Value = 0.4601490924567965
$ws.Range("Q6").Value = 41.846808053091
$ws.Range("R6").Value = 376.621272477819
$ws.Range("S6").Value = 0.03824632271828496
$ws.Range("T6").Value = 0.04151040203103293
$ws.Range("I7").Value = 0.2448858138641327
$ws.Range("J7").Value = 0.2525342032254661
$ws.Range("M7").Value = 0.4260053333333333
$ws.Range("N7").Value = 1.278016
$ws.Range("O7").Value = 0.02405532912416773
$ws.Range("P7").Value = 0.02531756756689831
$ws.Range("Q7").Value = 6.445365978147557
$ws.Range("R7").Value = 58.008293803328
$ws.Range("S7").Value = 0.005890808850341389
$ws.Range("T7").Value = 0.006393551753113567
$ws.Range("I8").Value = 0.2448858138641327
$ws.Range("J8").Value = 0.2525342032254661
$ws.Range("O8").Value = 0.05879323641880037
$ws.Range("P8").Value = 0.06187825274916518
$ws.Range("S8").Value = 0.0143976295501243
$ws.Range("T8").Value = 0.01562637525499444
$ws.Range("I9").Value = 0.2448858138641327
$ws.Range("J9").Value = 0.2525342032254661
$ws.Range("M9").Value = 5.850740666666667
$ws.Range("N9").Value = 17.552222
$ws.Range("O9").Value = 0.3303749538898241
$ws.Range("P9").Value = 0.3477104875323931
$ws.Range("Q9").Value = 88.52040547199178
$ws.Range("R9").Value = 796.683649247926
$ws.Range("S9").Value = 0.08090413946363491
$ws.Range("T9").Value = 0.08780879092213127
$ws.Range("I10").Value = 0.2448858138641327
$ws.Range("J10").Value = 0.2525342032254661
$ws.Range("M10").Value = 2.648771
$ws.Range("N10").Value = 5.297542
$ws.Range("O10").Value = 0.1495686865725097
$ws.Range("P10").Value = 0.1049445996947469
$ws.Range("Q10").Value = 40.07531631991434
$ws.Range("R10").Value = 240.451897919486
$ws.Range("S10").Value = 0.03662724953989843
$ws.Range("T10").Value = 0.02650210086672839
$ws.Range("I11").Value = 0.2448858138641327
$ws.Range("J11").Value = 0.2525342032254661
$ws.Range("M11").Value = 7.742685666666667
$ws.Range("N11").Value = 23.228057
$ws.Range("O11").Value = 0.4372077939946981
$ws.Range("P11").Value = 0.4601490924567965
$ws.Range("Q11").Value = 117.1451126795534
$ws.Range("R11").Value = 1054.306014115981
$ws.Range("S11").Value = 0.1070659864601337
$ws.Range("T11").Value = 0.1162033844284984
$ws.Range("G12").Value = 17.564497
$ws.Range("H12").Value = 52.69349099999999
$ws.Range("I12").Value = 0.2842934138044583
$ws.Range("J12").Value = 0.2931725999334087
$ws.Range("M12").Value = 0.4260053333333333
$ws.Range("N12").Value = 1.278016
$ws.Range("O12").Value = 0.02405532912416773
$ws.Range("P12").Value = 0.02531756756689831
$ws.Range("Q12").Value = 7.482569399317334
$ws.Range("R12").Value = 67.34312459385599
$ws.Range("S12").Value = 0.006838771636899452
$ws.Range("T12").Value = 0.007422417107577322
$ws.Range("G13").Value = 17.564497
$ws.Range("H13").Value = 52.69349099999999
$ws.Range("I13").Value = 0.2842934138044583
$ws.Range("J13").Value = 0.2931725999334087
$ws.Range("O13").Value = 0.05879323641880037
$ws.Range("P13").Value = 0.06187825274916518
$ws.Range("Q13").Value = 18.28802547008867
$ws.Range("R13").Value = 164.592229230798
$ws.Range("S13").Value = 0.01671452989011336
$ws.Range("T13").Value = 0.01814100823780935
$ws.Range("G14").Value = 17.564497
$ws.Range("H14").Value = 52.69349099999999
$ws.Range("I14").Value = 0.2842934138044583
$ws.Range("J14").Value = 0.2931725999334087
$ws.Range("M14").Value = 5.850740666666667
$ws.Range("N14").Value = 17.552222
$ws.Range("O14").Value = 0.3303749538898241
$ws.Range("P14").Value = 0.3477104875323931
$ws.Range("Q14").Value = 102.7653168874447
$ws.Range("R14").Value = 924.8878519870019
$ws.Range("S14").Value = 0.0939234234768286
$ws.Range("T14").Value = 0.1019391876539848
$ws.Range("G15").Value = 17.564497
$ws.Range("H15").Value = 52.69349099999999
$ws.Range("I15").Value = 0.2842934138044583
$ws.Range("J15").Value = 0.2931725999334087
$ws.Range("M15").Value = 2.648771
$ws.Range("N15").Value = 5.297542
$ws.Range("O15").Value = 0.1495686865725097
$ws.Range("P15").Value = 0.1049445996947469
$ws.Range("Q15").Value = 46.524330283187
$ws.Range("R15").Value = 279.145981699122
$ws.Range("S15").Value = 0.04252139250394783
$ws.Range("T15").Value = 0.03076688114147975
$ws.Range("G16").Value = 17.564497
$ws.Range("H16").Value = 52.69349099999999
$ws.Range("I16").Value = 0.2842934138044583
$ws.Range("J16").Value = 0.2931725999334087
$ws.Range("M16").Value = 7.742685666666667
$ws.Range("N16").Value = 23.228057
$ws.Range("O16").Value = 0.4372077939946981
$ws.Range("P16").Value = 0.4601490924567965
$ws.Range("Q16").Value = 135.9963791641097
$ws.Range("R16").Value = 1223.967412476987
$ws.Range("S16").Value = 0.124295296296669
$ws.Range("T16").Value = 0.1349031057925575
$ws.Range("G17").Value = 5.613580499999999
$ws.Range("H17").Value = 11.227161
$ws.Range("I17").Value = 0.09085964511315853
$ws.Range("J17").Value = 0.06246494429911598
$ws.Range("M17").Value = 0.4260053333333333
$ws.Range("N17").Value = 1.278016
$ws.Range("O17").Value = 0.02405532912416773
$ws.Range("P17").Value = 0.02531756756689831
$ws.Range("Q17").Value = 2.391415232096
$ws.Range("R17").Value = 14.348491392576
$ws.Range("S17").Value = 0.002185658667302106
$ws.Range("T17").Value = 0.001581460447855408
$ws.Range("G18").Value = 5.613580499999999
$ws.Range("H18").Value = 11.227161
$ws.Range("I18").Value = 0.09085964511315853
$ws.Range("J18").Value = 0.06246494429911598
$ws.Range("O18").Value = 0.05879323641880037
$ws.Range("P18").Value = 0.06187825274916518
$ws.Range("Q18").Value = 5.844818850343
$ws.Range("R18").Value = 35.068913102058
$ws.Range("S18").Value = 0.005341932596066229
$ws.Range("T18").Value = 0.003865221611303224
$ws.Range("G19").Value = 5.613580499999999
$ws.Range("H19").Value = 11.227161
$ws.Range("I19").Value = 0.09085964511315853
$ws.Range("J19").Value = 0.06246494429911598
$ws.Range("M19").Value = 5.850740666666667
$ws.Range("N19").Value = 17.552222
$ws.Range("O19").Value = 0.3303749538898241
$ws.Range("P19").Value = 0.3477104875323931
$ws.Range("Q19").Value = 32.843603716957
$ws.Range("R19").Value = 197.061622301742
$ws.Range("S19").Value = 0.03001775106470553
$ws.Range("T19").Value = 0.0217197162359294
$ws.Range("G20").Value = 5.613580499999999
$ws.Range("H20").Value = 11.227161
$ws.Range("I20").Value = 0.09085964511315853
$ws.Range("J20").Value = 0.06246494429911598
$ws.Range("M20").Value = 2.648771
$ws.Range("N20").Value = 5.297542
$ws.Range("O20").Value = 0.1495686865725097
$ws.Range("P20").Value = 0.1049445996947469
$ws.Range("Q20").Value = 14.8690892345655
$ws.Range("R20").Value = 59.47635693826199
$ws.Range("S20").Value = 0.01358975778201947
$ws.Range("T20").Value = 0.006555358574425386
$ws.Range("G21").Value = 5.613580499999999
$ws.Range("H21").Value = 11.227161
$ws.Range("I21").Value = 0.09085964511315853
$ws.Range("J21").Value = 0.06246494429911598
$ws.Range("M21").Value = 7.742685666666667
$ws.Range("N21").Value = 23.228057
$ws.Range("O21").Value = 0.4372077939946981
$ws.Range("P21").Value = 0.4601490924567965
$ws.Range("Q21").Value = 43.46418927602949
$ws.Range("R21").Value = 260.785135656177
$ws.Range("S21").Value = 0.03972454500306519
$ws.Range("T21").Value = 0.02874318742960256
$ws.Range("G22").Value = 18.07044533333334
$ws.Range("H22").Value = 54.211336
$ws.Range("I22").Value = 0.2924825341016128
$ws.Range("J22").Value = 0.3016174867021735
$ws.Range("M22").Value = 0.4260053333333333
$ws.Range("N22").Value = 1.278016
$ws.Range("O22").Value = 0.02405532912416773
$ws.Range("P22").Value = 0.02531756756689831
$ws.Range("Q22").Value = 7.698106087708446
$ws.Range("R22").Value = 69.28295478937601
$ws.Range("S22").Value = 0.007035763620884906
$ws.Range("T22").Value = 0.00763622109894033
$ws.Range("G23").Value = 18.07044533333334
$ws.Range("H23").Value = 54.211336
$ws.Range("I23").Value = 0.2924825341016128
$ws.Range("J23").Value = 0.3016174867021735
$ws.Range("O23").Value = 0.05879323641880037
$ws.Range("P23").Value = 0.06187825274916518
$ws.Range("Q23").Value = 18.81481516446756
$ws.Range("R23").Value = 169.333336480208
$ws.Range("S23").Value = 0.01719599477580596
$ws.Range("T23").Value = 0.01866356307572506
$ws.Range("G24").Value = 18.07044533333334
$ws.Range("H24").Value = 54.211336
$ws.Range("I24").Value = 0.2924825341016128
$ws.Range("J24").Value = 0.3016174867021735
$ws.Range("M24").Value = 5.850740666666667
$ws.Range("N24").Value = 17.552222
$ws.Range("O24").Value = 0.3303749538898241
$ws.Range("P24").Value = 0.3477104875323931
$ws.Range("Q24").Value = 105.7254893765102
$ws.Range("R24").Value = 951.5294043885921
$ws.Range("S24").Value = 0.09662890371739924
$ws.Range("T24").Value = 0.1048755633495079
$ws.Range("G25").Value = 18.07044533333334
$ws.Range("H25").Value = 54.211336
$ws.Range("I25").Value = 0.2924825341016128
$ws.Range("J25").Value = 0.3016174867021735
$ws.Range("M25").Value = 2.648771
$ws.Range("N25").Value = 5.297542
$ws.Range("O25").Value = 0.1495686865725097
$ws.Range("P25").Value = 0.1049445996947469
$ws.Range("Q25").Value = 47.86447155601867
$ws.Range("R25").Value = 287.186829336112
$ws.Range("S25").Value = 0.04374622847097751
$ws.Range("T25").Value = 0.03165312640289523
$ws.Range("G26").Value = 18.07044533333334
$ws.Range("H26").Value = 54.211336
$ws.Range("I26").Value = 0.2924825341016128
$ws.Range("J26").Value = 0.3016174867021735
$ws.Range("M26").Value = 7.742685666666667
$ws.Range("N26").Value = 23.228057
$ws.Range("O26").Value = 0.4372077939946981
$ws.Range("P26").Value = 0.4601490924567965
$ws.Range("Q26").Value = 139.9137780726836
$ws.Range("R26").Value = 1259.224002654152
$ws.Range("S26").Value = 0.1278756435165452
$ws.Range("T26").Value = 0.138789012775105

Write-Output "Updated 278 cells with new TPM values"
